# The commit deletes the slide titled "ARL and EGT Co-Adaptation and Influence"
# (internal p:sldId = 278) from the deck. All remaining slides keep their
# content unchanged; they simply shift up by one position once this slide
# is removed.

$p = $ppt.ActivePresentation

$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    if ($s.SlideID -eq 278) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    # Fallback: locate by title text in case SlideID lookup is unavailable.
    for ($i = 1; $i -le $p.Slides.Count; $i++) {
        $s = $p.Slides.Item($i)
        try {
            $t = $s.Shapes.Title.TextFrame.TextRange.Text
        } catch {
            $t = ""
        }
        if ($t -eq "ARL and EGT Co-Adaptation and Influence") {
            $targetIndex = $i
            break
        }
    }
}

if ($targetIndex -ne -1) {
    $p.Slides.Item($targetIndex).Delete()
}
